$wb = $excel.ActiveWorkbook

# --- Update data values in "Method" sheet (row 11 and 12) ---
$wsMethod = $wb.Worksheets.Item("Method")
$wsMethod.Range("F11").Value = 1
$wsMethod.Range("G11").Value = 2
$wsMethod.Range("H11").Value = 4
$wsMethod.Range("I11").Value = 8
$wsMethod.Range("J11").Value = 16
$wsMethod.Range("K11").Value = 32

$wsMethod.Range("G12").Value = 2
$wsMethod.Range("H12").Value = 4
$wsMethod.Range("I12").Value = 8
$wsMethod.Range("J12").Value = 16
$wsMethod.Range("K12").Value = 32

# --- Update selections / active sheet ---
$wsMethod.Range("H17").Select()
$wsMethod.Activate()
